$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style from an untouched, default-styled cell, used to keep
# text-forced cells on the original (unstyled) cell format after the edit.
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.819.57'
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = '  -0.70%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.296.78'
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = '  -1.02%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.17'
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  -1.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.21'
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  -2.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.513'
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = '  +1.09%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.505'
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = '  -2.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.88'
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = '  -0.57%  '

$ws.Range("E11").Value = '  -0.52%  '

$ws.Range("E12").Value = '  +0.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.68'
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = '  -0.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = '  -2.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.644.88'
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = '  -1.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.297.99'
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.776'
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = '  -2.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.791.28'
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = '  -0.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.51'
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '  -5.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0906'
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = '  -0.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.05'
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '  -2.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.88'
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = '  -0.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.08'
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = '  +0.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = '  -1.38%  '

$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.02'
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = '  -0.30%  '

$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.41'
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = '  -1.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.14'
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = '  -1.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.90'
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = '  -1.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.00'
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = '  -2.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.73'
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = '  -4.09%  '

$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("E34").Value = '  -4.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.98'
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '  -3.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.09'
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  -4.99%  '

$ws.Range("E37").Value = '  -0.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0686'
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = '  -1.79%  '

$ws.Range("E39").Value = '  -2.38%  '

$ws.Range("E40").Value = '  -4.01%  '

$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("E42").Value = '  -0.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.014.82'
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = '  +1.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0283'
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = '  -1.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.09'
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '  -1.21%  '

$ws.Range("E46").Value = '  -4.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.07'
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = '  -3.10%  '

$ws.Range("E48").Value = '  -3.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.515.17'
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = '  -1.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.13'
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = '  -3.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.79'
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = '  -7.77%  '
